$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for every data row (2-130)
# from 2023-09-23 (45192) to 2023-10-03 (45202).
for ($r = 2; $r -le 130; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45192) {
        $cell.Value2 = 45202
    }
}

# Row 5 (A 67432-2021) gained a new observed species "Skogsknipprot",
# which bumps the Fridlysta (protected), Signalarter (signal species)
# and Alla arter (all species) counts by one each, and the species list
# in column R needs the new name inserted alphabetically-ish right
# after "Nästrot" and before "Svavelriska".
$ws.Cells.Item(5, 8).Value2 = 4
$ws.Cells.Item(5, 9).Value2 = 8
$ws.Cells.Item(5, 17).Value2 = 10

$r5 = $ws.Cells.Item(5, 18)
$oldList = $r5.Value2
$newList = $oldList.Replace("Nästrot`r`nSvavelriska", "Nästrot`r`nSkogsknipprot`r`nSvavelriska")
$r5.Value2 = $newList
